{"js": "// Rewrite the bullet list under \"KEY ACHIEVEMENTS AND IMPACT\" / \"Impact\" from\n// six job-duty-style bullets into four impact-focused accomplishment\n// statements, per the commit:\n//   \"Fix Key Achievements to use proper accomplishment statements\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading, then its \"Impact\"\n// sub-heading right after it. The six bullet paragraphs we need to rewrite\n// are the ones that immediately follow \"Impact\" (some of this section's\n// bullet text is duplicated verbatim elsewhere in the resume - e.g. in the\n// \"Partner - Siege Analytics\" experience entry - so we anchor on the\n// section heading instead of doing a bare text search).\nconst items = paragraphs.items;\nlet sectionHeadingIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    sectionHeadingIdx = i;\n    break;\n  }\n}\nif (sectionHeadingIdx === -1) {\n  throw new Error('Could not find \"KEY ACHIEVEMENTS AND IMPACT\" heading');\n}\n\nlet impactHeadingIdx = -1;\nfor (let i = sectionHeadingIdx + 1; i < items.length; i++) {\n  if (items[i].text.trim() === \"Impact\") {\n    impactHeadingIdx = i;\n    break;\n  }\n}\nif (impactHeadingIdx === -1) {\n  throw new Error('Could not find \"Impact\" sub-heading under Key Achievements');\n}\n\n// The old bullet block (6 paragraphs) runs from just after \"Impact\" up to\n// (but not including) the next heading-styled paragraph (\"TECHNICAL SKILLS\").\nconst oldBulletStart = impactHeadingIdx + 1;\nconst oldBullets = [];\nfor (let i = oldBulletStart; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (t.startsWith(\"\u2022\")) {\n    oldBullets.push(i);\n  } else {\n    break;\n  }\n}\n\nconst newBullets = [\n  \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  \"\u2022 $4.7M savings enabled nonprofit access\",\n  \"\u2022 Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\",\n  \"\u2022 Real-time collaboration at national scale\",\n];\n\n// Rewrite the text of the first len(newBullets) paragraphs in place (keeps\n// their paragraph formatting/run formatting intact), then delete any extra\n// leftover bullet paragraphs beyond that.\nfor (let i = 0; i < newBullets.length; i++) {\n  items[oldBullets[i]].insertText(newBullets[i], \"Replace\");\n}\nfor (let i = newBullets.length; i < oldBullets.length; i++) {\n  items[oldBullets[i]].delete();\n}\n\nawait context.sync();\n", "ps1": "# Rewrite the bullet list under \"KEY ACHIEVEMENTS AND IMPACT\" / \"Impact\" from\n# six job-duty-style bullets into four impact-focused accomplishment\n# statements, per the commit:\n#   \"Fix Key Achievements to use proper accomplishment statements\"\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n# Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading, then its \"Impact\"\n# sub-heading right after it. Some of this section's bullet text is\n# duplicated verbatim elsewhere in the resume (e.g. in the \"Partner - Siege\n# Analytics\" experience entry), so anchor on the section heading rather than\n# doing a bare text search/replace.\n$sectionHeadingIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $paras.Item($i).Range.Text.Trim()\n    if ($t -eq \"KEY ACHIEVEMENTS AND IMPACT\") {\n        $sectionHeadingIdx = $i\n        break\n    }\n}\nif ($sectionHeadingIdx -eq -1) {\n    throw \"Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading\"\n}\n\n$impactHeadingIdx = -1\nfor ($i = $sectionHeadingIdx + 1; $i -le $count; $i++) {\n    $t = $paras.Item($i).Range.Text.Trim()\n    if ($t -eq \"Impact\") {\n        $impactHeadingIdx = $i\n        break\n    }\n}\nif ($impactHeadingIdx -eq -1) {\n    throw \"Could not find 'Impact' sub-heading under Key Achievements\"\n}\n\n# The old bullet block runs from just after \"Impact\" up to (but not\n# including) the next non-bullet paragraph (\"TECHNICAL SKILLS\").\n$oldBulletIdxs = New-Object System.Collections.ArrayList\nfor ($i = $impactHeadingIdx + 1; $i -le $count; $i++) {\n    $t = $paras.Item($i).Range.Text.Trim()\n    if ($t.StartsWith(\"\u2022\")) {\n        [void]$oldBulletIdxs.Add($i)\n    } else {\n        break\n    }\n}\n\n$newBullets = @(\n    \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n    \"\u2022 `$4.7M savings enabled nonprofit access\",\n    \"\u2022 Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\",\n    \"\u2022 Real-time collaboration at national scale\"\n)\n\n# Rewrite the text of the first N paragraphs in place (keeps their paragraph\n# formatting/run formatting intact), then delete any extra leftover bullet\n# paragraphs beyond that. Delete from the end backwards so earlier indices\n# in $oldBulletIdxs stay valid.\nfor ($i = 0; $i -lt $newBullets.Count; $i++) {\n    $idx = $oldBulletIdxs[$i]\n    $paras.Item($idx).Range.Text = $newBullets[$i]\n}\nfor ($i = $oldBulletIdxs.Count - 1; $i -ge $newBullets.Count; $i--) {\n    $idx = $oldBulletIdxs[$i]\n    $paras.Item($idx).Range.Delete()\n}\n"}
